$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column I, remove column J header
$ws.Range("I1").Value = "roc_auc_overall_validation"
$ws.Range("J1:J3").Clear()

# Row 2 becomes eval001 (new row)
$ws.Range("A2").Value = "11:19AM CET on Nov 30, 2022"
$ws.Range("B2").Value = "eval001"
$ws.Range("C2").Value = "/home/sascha/Projects/BPDP-Dataset_2022/BPDP_Dataset/"
$ws.Range("D2").Value = "GradientBoostingClassifier"
$ws.Range("E2").Value = 0.9801980198019802
$ws.Range("F2").Value = 0.853
$ws.Range("G2").Value = 0.06
$ws.Range("H2").Value = "/home/sascha/Projects/favel/Evaluation/eval001"
$ws.Range("I2").Value = 0.6747572815533981

# Row 3 is the updated eval002 entry
$ws.Range("A3").Value = "11:19AM CET on Nov 30, 2022"
$ws.Range("B3").Value = "eval002"
$ws.Range("C3").Value = "/home/sascha/Projects/BPDP-Dataset_2022/BPDP_Dataset/"
$ws.Range("D3").Value = "RandomForestClassifier"
$ws.Range("E3").Value = 0.8811881188118812
$ws.Range("F3").Value = 0.9268571428571428
$ws.Range("G3").Value = 0.05
$ws.Range("H3").Value = "/home/sascha/Projects/favel/Evaluation/eval002"
$ws.Range("I3").Value = 0.6699029126213593
